{"js": "// The document contains a date heading and a grid of two-digit by\n// two-digit multiplication prompts. Each source run's text is replaced\n// with the text from the commit's target revision, matched 1:1 in\n// document order (every source string below is unique in the doc, so a\n// simple exact-text search/replace is unambiguous).\nconst replacements = [\n  [\"2026-02-19 Thursday\", \"2026-02-20 Friday\"],\n  [\"41\u00d762=\", \"14\u00d761=\"],\n  [\"56\u00d716=\", \"13\u00d748=\"],\n  [\"51\u00d719=\", \"18\u00d774=\"],\n  [\"19\u00d728=\", \"70\u00d731=\"],\n  [\"93\u00d763=\", \"48\u00d743=\"],\n  [\"81\u00d793=\", \"40\u00d780=\"],\n  [\"95\u00d740=\", \"29\u00d747=\"],\n  [\"62\u00d793=\", \"85\u00d792=\"],\n  [\"49\u00d713=\", \"79\u00d780=\"],\n  [\"77\u00d743=\", \"59\u00d764=\"],\n  [\"49\u00d779=\", \"29\u00d754=\"],\n  [\"15\u00d794=\", \"15\u00d742=\"],\n  [\"85\u00d715=\", \"89\u00d733=\"],\n  [\"43\u00d783=\", \"20\u00d759=\"],\n  [\"94\u00d720=\", \"84\u00d786=\"],\n  [\"13\u00d782=\", \"55\u00d734=\"],\n  [\"81\u00d744=\", \"72\u00d757=\"],\n  [\"27\u00d738=\", \"42\u00d799=\"],\n  [\"90\u00d796=\", \"94\u00d731=\"],\n  [\"40\u00d799=\", \"73\u00d715=\"],\n  [\"78\u00d752=\", \"52\u00d743=\"],\n  [\"95\u00d781=\", \"30\u00d746=\"],\n  [\"39\u00d763=\", \"77\u00d756=\"],\n  [\"93\u00d748=\", \"48\u00d764=\"],\n  [\"54\u00d799=\", \"71\u00d719=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document is a date heading followed by a 5x25 grid of two-digit by\n# two-digit multiplication prompts (laid out as five rows of text separated\n# by blank rows). Each source string below is unique within the document,\n# so an exact-text Find/Replace (ReplaceAll) for every (old,new) pair is\n# unambiguous and reproduces the whole diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-02-19 Thursday\", \"2026-02-20 Friday\"),\n    @(\"41\u00d762=\", \"14\u00d761=\"),\n    @(\"56\u00d716=\", \"13\u00d748=\"),\n    @(\"51\u00d719=\", \"18\u00d774=\"),\n    @(\"19\u00d728=\", \"70\u00d731=\"),\n    @(\"93\u00d763=\", \"48\u00d743=\"),\n    @(\"81\u00d793=\", \"40\u00d780=\"),\n    @(\"95\u00d740=\", \"29\u00d747=\"),\n    @(\"62\u00d793=\", \"85\u00d792=\"),\n    @(\"49\u00d713=\", \"79\u00d780=\"),\n    @(\"77\u00d743=\", \"59\u00d764=\"),\n    @(\"49\u00d779=\", \"29\u00d754=\"),\n    @(\"15\u00d794=\", \"15\u00d742=\"),\n    @(\"85\u00d715=\", \"89\u00d733=\"),\n    @(\"43\u00d783=\", \"20\u00d759=\"),\n    @(\"94\u00d720=\", \"84\u00d786=\"),\n    @(\"13\u00d782=\", \"55\u00d734=\"),\n    @(\"81\u00d744=\", \"72\u00d757=\"),\n    @(\"27\u00d738=\", \"42\u00d799=\"),\n    @(\"90\u00d796=\", \"94\u00d731=\"),\n    @(\"40\u00d799=\", \"73\u00d715=\"),\n    @(\"78\u00d752=\", \"52\u00d743=\"),\n    @(\"95\u00d781=\", \"30\u00d746=\"),\n    @(\"39\u00d763=\", \"77\u00d756=\"),\n    @(\"93\u00d748=\", \"48\u00d764=\"),\n    @(\"54\u00d799=\", \"71\u00d719=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace) -- Wrap:=wdFindContinue(1), Replace:=wdReplaceAll(2)\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n"}
